# Reproduce the commit: reorder metadata rows 15-17 (cyclic shift) and
# update the active selection on the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rotate the attribute rows 15..17 -------------------------------------
# Before:  A15/B15 = t_for / "...Foragin..."
#          A16/B16 = n_def / "...Defecation..."
#          A17/B17 = t_beb / "...Drinking..."
# After:   A15/B15 = n_def / "...Defecation..."
#          A16/B16 = t_beb / "...Drinking..."
#          A17/B17 = t_for / "...Foragin..."
$a15 = $ws.Range("A15").Value()
$b15 = $ws.Range("B15").Value()
$a16 = $ws.Range("A16").Value()
$b16 = $ws.Range("B16").Value()
$a17 = $ws.Range("A17").Value()
$b17 = $ws.Range("B17").Value()

$ws.Range("A15").Value = $a16
$ws.Range("B15").Value = $b16
$ws.Range("A16").Value = $a17
$ws.Range("B16").Value = $b17
$ws.Range("A17").Value = $a15
$ws.Range("B17").Value = $b15

# --- Update the worksheet view / selection --------------------------------
# Move the active selection to a single cell (A21) and let the view's
# topLeftCell reset naturally (no longer pinned to A5).
$ws.Range("A21").Select()
